$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update "Salario Basico" values (column G, rows 16-18) from 1500000 to 1800000
$ws.Range("G16").Value = 1800000
$ws.Range("G17").Value = 1800000
$ws.Range("G18").Value = 1800000

# Update "Periodo Mora" values (column E) - swap order so E16 shows 2502 and E18 shows 2504
$ws.Range("E16").Value = "2502"
$ws.Range("E18").Value = "2504"
